$wb = $excel.ActiveWorkbook

$id = "2593015e-e604-4d8a-a976-6453ef59a653"
$mdName = "$id.md"
$zhXlf = "$id.75fc1b361e66756fecd256450813964da5ae816f.zh-cn.xlf"
$deXlf = "$id.75fc1b361e66756fecd256450813964da5ae816f.de-de.xlf"

$mdAddr = "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/$mdName"
$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b1cb6875ac11b04c2ae16426321b369a9d2d2a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbc483f6f836dd93f3633908624a49a71bfb0286/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ---------------------------------------------------------------
# Sheet "Overview": add row 4
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Add($ov.Range("A4"), $mdAddr, "", "", $mdName)
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-48-13 14:48:14"

# ---------------------------------------------------------------
# Sheet "zh-cn": add row 4
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Add($zh.Range("A4"), $mdAddr, "", "", $mdName)
$zh.Hyperlinks.Add($zh.Range("B4"), $mdAddr, "", "", ".md")
$zh.Range("C4").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D4"), $zhXlfAddr, "", "", $zhXlf)
$zh.Range("E4").Value = "2016-03-13 14:48:11"
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "Include"

# ---------------------------------------------------------------
# Sheet "de-de": add row 4
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Add($de.Range("A4"), $mdAddr, "", "", $mdName)
$de.Hyperlinks.Add($de.Range("B4"), $mdAddr, "", "", ".md")
$de.Range("C4").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D4"), $deXlfAddr, "", "", $deXlf)
$de.Range("E4").Value = "2016-03-13 14:48:14"
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "Include"

Write-Host "Row 4 added to Overview, zh-cn and de-de sheets"
